$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "Term Opening Hours (...)" day headers to just the short
# weekday/PH label ("Mon", "Tue", ... "PH").
$ws.Range("F1").Value = "Mon"
$ws.Range("G1").Value = "Tue"
$ws.Range("H1").Value = "Wed"
$ws.Range("I1").Value = "Thu"
$ws.Range("J1").Value = "Fri"
$ws.Range("K1").Value = "Sat"
$ws.Range("L1").Value = "Sun"
$ws.Range("M1").Value = "PH"

# Leave the cursor where the author left it after editing.
$ws.Range("G14").Select()
